$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$c = $ws1.Cells.Item(4,136)
$c.Value = "03/01/2023 "
Write-Host "[$($c.Value2)]"
